$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the "Name" label in A5 (reuses existing shared string "Name")
$ws.Range("A5").Value = "Name"

# Update B5 from "justatest" to "name"
$ws.Range("B5").Value = "name"
